# Update probability matrix values on Sheet1 to reflect the refreshed
# simulation results (more games simulated, updated transition probabilities).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1849056603773585
$ws.Range("C2").Value = 0.5584905660377358
$ws.Range("J2").Value = 0.03018867924528302
$ws.Range("P2").Value = 0.1132075471698113
$ws.Range("S2").Value = 0.1132075471698113
$ws.Range("B3").Value = 0.006711409395973154
$ws.Range("C3").Value = 0.01342281879194631
$ws.Range("J3").Value = 0.01342281879194631
$ws.Range("P3").Value = 0.7181208053691275
$ws.Range("S3").Value = 0.2483221476510067
$ws.Range("J4").Value = 0.03125
$ws.Range("P4").Value = 0.65625
$ws.Range("S4").Value = 0.3125
$ws.Range("B6").Value = 0.09246575342465753
$ws.Range("D6").Value = 0.003424657534246575
$ws.Range("E6").Value = 0.003424657534246575
$ws.Range("F6").Value = 0.08561643835616438
$ws.Range("J6").Value = 0.2671232876712329
$ws.Range("O6").Value = 0.0273972602739726
$ws.Range("Q6").Value = 0.1678082191780822
$ws.Range("R6").Value = 0.03767123287671233
$ws.Range("S6").Value = 0.3150684931506849
$ws.Range("B7").Value = 0.09322033898305085
$ws.Range("D7").Value = 0.01694915254237288
$ws.Range("F7").Value = 0.0635593220338983
$ws.Range("J7").Value = 0.1016949152542373
$ws.Range("O7").Value = 0.02966101694915254
$ws.Range("Q7").Value = 0.1779661016949153
$ws.Range("R7").Value = 0.0847457627118644
$ws.Range("S7").Value = 0.4322033898305085
$ws.Range("B8").Value = 0.07721280602636535
$ws.Range("D8").Value = 0.01129943502824859
$ws.Range("F8").Value = 0.06591337099811675
$ws.Range("J8").Value = 0.0847457627118644
$ws.Range("O8").Value = 0.04143126177024482
$ws.Range("Q8").Value = 0.1713747645951036
$ws.Range("R8").Value = 0.08097928436911488
$ws.Range("S8").Value = 0.4670433145009416
$ws.Range("B9").Value = 0.07352941176470588
$ws.Range("D9").Value = 0.02941176470588235
$ws.Range("F9").Value = 0.07352941176470588
$ws.Range("J9").Value = 0.125
$ws.Range("O9").Value = 0.04411764705882353
$ws.Range("Q9").Value = 0.1764705882352941
$ws.Range("R9").Value = 0.1102941176470588
$ws.Range("S9").Value = 0.3676470588235294
$ws.Range("B10").Value = 0.09416666666666666
$ws.Range("D10").Value = 0.01583333333333333
$ws.Range("E10").Value = 0.001666666666666667
$ws.Range("F10").Value = 0.08166666666666667
$ws.Range("J10").Value = 0.1033333333333333
$ws.Range("O10").Value = 0.01416666666666667
$ws.Range("Q10").Value = 0.2283333333333333
$ws.Range("R10").Value = 0.06666666666666667
$ws.Range("S10").Value = 0.3941666666666667
$ws.Range("G11").Value = 0.15
$ws.Range("J11").Value = 0.1078947368421053
$ws.Range("K11").Value = 0.2263157894736842
$ws.Range("L11").Value = 0.4973684210526316
$ws.Range("S11").Value = 0.01842105263157895
$ws.Range("G12").Value = 0.7512953367875648
$ws.Range("J12").Value = 0.1658031088082902
$ws.Range("K12").Value = 0.02072538860103627
$ws.Range("L12").Value = 0.02072538860103627
$ws.Range("S12").Value = 0.04145077720207254
$ws.Range("G13").Value = 0.6779661016949152
$ws.Range("J13").Value = 0.288135593220339
$ws.Range("S13").Value = 0.03389830508474576
$ws.Range("F15").Value = 0.024
$ws.Range("H15").Value = 0.216
$ws.Range("I15").Value = 0.024
$ws.Range("J15").Value = 0.316
$ws.Range("K15").Value = 0.08
$ws.Range("M15").Value = 0.024
$ws.Range("O15").Value = 0.104
$ws.Range("S15").Value = 0.212
$ws.Range("F16").Value = 0.02564102564102564
$ws.Range("H16").Value = 0.2371794871794872
$ws.Range("I16").Value = 0.08333333333333333
$ws.Range("J16").Value = 0.3461538461538461
$ws.Range("K16").Value = 0.1153846153846154
$ws.Range("M16").Value = 0.01923076923076923
$ws.Range("O16").Value = 0.04487179487179487
$ws.Range("S16").Value = 0.1282051282051282
$ws.Range("F17").Value = 0.02536997885835095
$ws.Range("H17").Value = 0.200845665961945
$ws.Range("I17").Value = 0.06342494714587738
$ws.Range("J17").Value = 0.3974630021141649
$ws.Range("K17").Value = 0.1141649048625793
$ws.Range("M17").Value = 0.02959830866807611
$ws.Range("O17").Value = 0.07399577167019028
$ws.Range("S17").Value = 0.09513742071881606
$ws.Range("F18").Value = 0.05389221556886228
$ws.Range("H18").Value = 0.2035928143712575
$ws.Range("I18").Value = 0.05988023952095808
$ws.Range("J18").Value = 0.3832335329341318
$ws.Range("K18").Value = 0.125748502994012
$ws.Range("M18").Value = 0.02994011976047904
$ws.Range("O18").Value = 0.05988023952095808
$ws.Range("S18").Value = 0.08383233532934131
$ws.Range("F19").Value = 0.0400604686318972
$ws.Range("H19").Value = 0.236583522297808
$ws.Range("I19").Value = 0.0589569160997732
$ws.Range("J19").Value = 0.3371126228269085
$ws.Range("K19").Value = 0.1307634164777022
$ws.Range("M19").Value = 0.02343159486016629
$ws.Range("N19").Value = 0.0007558578987150416
$ws.Range("O19").Value = 0.06500377928949358
$ws.Range("S19").Value = 0.1073318216175359
